$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need an explicit
# text format first, otherwise Excel auto-converts the assigned
# string into a numeric value (losing formatting like trailing zeros).
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '26.588.31'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.820.56'
$ws.Range('E3').Value = '  +1.48%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = '1.009'
$ws.Range('D6').Value = '305.66'
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('D7').Value = '0.4674'
$ws.Range('E7').Value = '  +2.34%  '
$ws.Range('D8').Value = '0.3589'
$ws.Range('E8').Value = '  -1.03%  '
$ws.Range('D9').Value = '0.07124'
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('D10').Value = '0.9017'
$ws.Range('E10').Value = '  +2.66%  '
$ws.Range('D11').Value = '0.07801'
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').Value = '1.868.50'
$ws.Range('E13').Value = '  +4.15%  '
$ws.Range('D14').Value = '5.250'
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '87.25'
$ws.Range('E16').Value = '  +2.62%  '
$ws.Range('D17').Value = '1.011'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').Value = '0.000008547'
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').Value = '26.641.22'
$ws.Range('E20').Value = '  +0.69%  '
$ws.Range('D21').Value = '14.14'
$ws.Range('E21').Value = '  -0.78%  '
$ws.Range('D22').Value = '5.004'
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('D23').Value = '10.54'
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('D24').Value = '1.941'
$ws.Range('E24').Value = '  -2.03%  '
$ws.Range('D25').Value = '152.04'
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').Value = '17.88'
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  -3.38%  '
$ws.Range('D28').Value = '113.45'
$ws.Range('E28').Value = '  +1.35%  '
$ws.Range('D29').Value = '4.785'
$ws.Range('E29').Value = '  -1.25%  '
$ws.Range('D30').Value = '0.08801'
$ws.Range('E30').Value = '  +1.64%  '
$ws.Range('D31').Value = '3.148'
$ws.Range('E31').Value = '  +2.98%  '
$ws.Range('D32').Value = '2.753'
$ws.Range('E32').Value = '  +2.98%  '
$ws.Range('D33').Value = '0.7273'
$ws.Range('E33').Value = '  +0.48%  '
$ws.Range('D34').Value = '4.434'
$ws.Range('E34').Value = '  -0.26%  '
$ws.Range('E35').Value = '  +1.07%  '
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').Value = '0.01922'
$ws.Range('E37').Value = '  -1.01%  '
$ws.Range('D38').Value = '2.917'
$ws.Range('E38').Value = '  +1.82%  '
$ws.Range('D39').Value = '0.05106'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').Value = '6.834'
$ws.Range('E40').Value = '  -0.92%  '
$ws.Range('D41').Value = '0.5028'
$ws.Range('E41').Value = '  -3.54%  '
$ws.Range('D42').Value = '0.1493'
$ws.Range('E42').Value = '  -2.07%  '
$ws.Range('D43').Value = '7.976'
$ws.Range('E43').Value = '  -0.49%  '
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('D45').Value = '0.4649'
$ws.Range('E45').Value = '  -0.81%  '
$ws.Range('D46').Value = '10.01'
$ws.Range('D47').Value = '98.51'
$ws.Range('E47').Value = '  -1.38%  '
$ws.Range('D48').Value = '1.554'
$ws.Range('E48').Value = '  -1.99%  '
$ws.Range('D49').Value = '0.06000'
$ws.Range('E49').Value = '  +0.66%  '
$ws.Range('D50').Value = '63.57'
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('D51').Value = '35.64'
$ws.Range('E51').Value = '  -1.57%  '
